$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1476.814
$ws.Range("J17").Value = 1476.814
$ws.Range("L17").Value = 4430.442
$ws.Range("N17").Value = -4766.442
$ws.Range("H51").Value = 4965.684
$ws.Range("I51").Value = 4241.222
$ws.Range("J51").Value = 5617.7
$ws.Range("K51").Value = 4241.222
$ws.Range("L51").Value = 5617.7
$ws.Range("M51").Value = -3757.222
$ws.Range("N51").Value = -6585.7
$ws.Range("H62").Value = 3690.7273
$ws.Range("I62").Value = 3285.7144
$ws.Range("K62").Value = 3285.7144
$ws.Range("M62").Value = -2661.7144
$ws.Range("H65").Value = 3690.7273
$ws.Range("I65").Value = 3285.7144
$ws.Range("K65").Value = 16428.572
$ws.Range("M65").Value = -13308.572
$ws.Range("H86").Value = 2058300.9
$ws.Range("I86").Value = 3238861.8
$ws.Range("K86").Value = 3238861.8
$ws.Range("M86").Value = -3237738.8
$ws.Range("H89").Value = 2058300.9
$ws.Range("I89").Value = 3238861.8
$ws.Range("K89").Value = 16194309
$ws.Range("M89").Value = -16188693

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10453.789
$ws.Range("I32").Value = 10433.427
$ws.Range("J32").Value = 10555.6
$ws.Range("K32").Value = 10433.427
$ws.Range("L32").Value = 10555.6
$ws.Range("M32").Value = -10146.427
$ws.Range("N32").Value = -11129.6
$ws.Range("H45").Value = 5457.846
$ws.Range("I45").Value = 5157.3335
$ws.Range("J45").Value = 5715.4287
$ws.Range("K45").Value = 5157.3335
$ws.Range("L45").Value = 5715.4287
$ws.Range("M45").Value = -4780.3335
$ws.Range("N45").Value = -6469.4287
$ws.Range("H102").Value = 687134
$ws.Range("I102").Value = 1055364.8
$ws.Range("K102").Value = 1055364.8
$ws.Range("M102").Value = -1053742.8
$ws.Range("H122").Value = 6497.923
$ws.Range("J122").Value = 8800
$ws.Range("L122").Value = 26400
$ws.Range("N122").Value = -31300
$ws.Range("H132").Value = 16790.674
$ws.Range("I132").Value = 28722.62
$ws.Range("K132").Value = 86167.86
$ws.Range("M132").Value = -83637.86

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2495
$ws.Range("I105").Value = 2295.6667
$ws.Range("K105").Value = 2295.6667
$ws.Range("M105").Value = -548.6667000000002
$ws.Range("H134").Value = 1327.4839
$ws.Range("I134").Value = 798.4483
$ws.Range("K134").Value = 2395.3449
$ws.Range("M134").Value = 139.6550999999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6167.1226
$ws.Range("I31").Value = 1559.1538
$ws.Range("J31").Value = 7831.1113
$ws.Range("K31").Value = 1559.1538
$ws.Range("L31").Value = 7831.1113
$ws.Range("M31").Value = -1264.1538
$ws.Range("N31").Value = -8421.1113
$ws.Range("H34").Value = 6167.1226
$ws.Range("I34").Value = 1559.1538
$ws.Range("J34").Value = 7831.1113
$ws.Range("K34").Value = 1559.1538
$ws.Range("L34").Value = 7831.1113
$ws.Range("M34").Value = -1357.1538
$ws.Range("N34").Value = -8235.1113
$ws.Range("H58").Value = 273648.44
$ws.Range("I58").Value = 501755.16
$ws.Range("J58").Value = 5287.5884
$ws.Range("K58").Value = 501755.16
$ws.Range("L58").Value = 5287.5884
$ws.Range("M58").Value = -501552.16
$ws.Range("N58").Value = -5693.5884
$ws.Range("H86").Value = 7196.478
$ws.Range("I86").Value = 6327.706
$ws.Range("K86").Value = 6327.706
$ws.Range("M86").Value = -5204.706
$ws.Range("H89").Value = 7196.478
$ws.Range("I89").Value = 6327.706
$ws.Range("K89").Value = 31638.53
$ws.Range("M89").Value = -26022.53
$ws.Range("H134").Value = 1900.25
$ws.Range("I134").Value = 1544.2188
$ws.Range("K134").Value = 4632.6564
$ws.Range("M134").Value = -2097.6564
$ws.Range("H136").Value = 273648.44
$ws.Range("I136").Value = 501755.16
$ws.Range("J136").Value = 5287.5884
$ws.Range("K136").Value = 1505265.48
$ws.Range("L136").Value = 15862.7652
$ws.Range("M136").Value = -1502715.48
$ws.Range("N136").Value = -20962.7652
$ws.Range("H141").Value = 79443.07000000001
$ws.Range("J141").Value = 83120.83
$ws.Range("L141").Value = 83120.83
$ws.Range("N141").Value = -93480.83

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 312.8
$ws.Range("J2").Value = 356.25
$ws.Range("L2").Value = 2137.5
$ws.Range("N2").Value = -2363.5
$ws.Range("H7").Value = 302.63635
$ws.Range("I7").Value = 310
$ws.Range("J7").Value = 283
$ws.Range("K7").Value = 930
$ws.Range("L7").Value = 849
$ws.Range("M7").Value = -818
$ws.Range("N7").Value = -1073
$ws.Range("H17").Value = 4414.6665
$ws.Range("J17").Value = 4414.6665
$ws.Range("L17").Value = 13243.9995
$ws.Range("N17").Value = -13581.9995
$ws.Range("H132").Value = 6203.364
$ws.Range("I132").Value = 1373.1666
$ws.Range("J132").Value = 11999.6
$ws.Range("K132").Value = 12358.4994
$ws.Range("L132").Value = 107996.4
$ws.Range("M132").Value = -9828.499400000001
$ws.Range("N132").Value = -113056.4
$ws.Range("H133").Value = 26496.084
$ws.Range("I133").Value = 9499.5
$ws.Range("K133").Value = 28498.5
$ws.Range("M133").Value = -23438.5
$ws.Range("H139").Value = 2764.4707
$ws.Range("I139").Value = 1899.8
$ws.Range("J139").Value = 3999.7144
$ws.Range("K139").Value = 5699.4
$ws.Range("L139").Value = 11999.1432
$ws.Range("M139").Value = -559.3999999999996
$ws.Range("N139").Value = -22279.1432

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 694045.9
$ws.Range("I122").Value = 2203379.5
$ws.Range("J122").Value = 7985.091
$ws.Range("K122").Value = 6610138.5
$ws.Range("L122").Value = 23955.273
$ws.Range("M122").Value = -6607688.5
$ws.Range("N122").Value = -28855.273
$ws.Range("H132").Value = 94680.69500000001
$ws.Range("I132").Value = 124721.234
$ws.Range("K132").Value = 374163.702
$ws.Range("M132").Value = -371633.702
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7798.7
$ws.Range("I7").Value = 3989
$ws.Range("J7").Value = 8222
$ws.Range("K7").Value = 3989
$ws.Range("L7").Value = 8222
$ws.Range("M7").Value = -3877
$ws.Range("N7").Value = -8446
$ws.Range("H61").Value = 4134.4546
$ws.Range("I61").Value = 3311.8572
$ws.Range("J61").Value = 5574
$ws.Range("K61").Value = 3311.8572
$ws.Range("L61").Value = 5574
$ws.Range("M61").Value = -3109.8572
$ws.Range("N61").Value = -5978
$ws.Range("H113").Value = 4134.4546
$ws.Range("I113").Value = 3311.8572
$ws.Range("J113").Value = 5574
$ws.Range("K113").Value = 3311.8572
$ws.Range("L113").Value = 5574
$ws.Range("M113").Value = -1141.8572
$ws.Range("N113").Value = -9914
$ws.Range("H126").Value = 7798.7
$ws.Range("I126").Value = 3989
$ws.Range("J126").Value = 8222
$ws.Range("K126").Value = 11967
$ws.Range("L126").Value = 24666
$ws.Range("M126").Value = -9497
$ws.Range("N126").Value = -29606
$ws.Range("H132").Value = 4568.64
$ws.Range("I132").Value = 4170.1055
$ws.Range("K132").Value = 12510.3165
$ws.Range("M132").Value = -9980.316499999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5861.5454
$ws.Range("I62").Value = 3925.2856
$ws.Range("J62").Value = 9250
$ws.Range("K62").Value = 3925.2856
$ws.Range("L62").Value = 9250
$ws.Range("M62").Value = -3301.2856
$ws.Range("N62").Value = -10498
$ws.Range("H65").Value = 5861.5454
$ws.Range("I65").Value = 3925.2856
$ws.Range("J65").Value = 9250
$ws.Range("K65").Value = 19626.428
$ws.Range("L65").Value = 46250
$ws.Range("M65").Value = -16506.428
$ws.Range("N65").Value = -52490
$ws.Range("H107").Value = 2569.3872
$ws.Range("I107").Value = 2669.2856
$ws.Range("J107").Value = 2359.6
$ws.Range("K107").Value = 8007.8568
$ws.Range("L107").Value = 7078.799999999999
$ws.Range("M107").Value = -6087.8568
$ws.Range("N107").Value = -10918.8
$ws.Range("H113").Value = 918.5625
$ws.Range("I113").Value = 902.5454999999999
$ws.Range("K113").Value = 2707.6365
$ws.Range("M113").Value = -537.6364999999996
$ws.Range("H122").Value = 5304.8
$ws.Range("I122").Value = 4821.731
$ws.Range("K122").Value = 14465.193
$ws.Range("M122").Value = -12015.193
$ws.Range("H126").Value = 4418.2
$ws.Range("I126").Value = 4240.7144
$ws.Range("K126").Value = 12722.1432
$ws.Range("M126").Value = -10252.1432
